$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update revised historical values (rows 218-220: columns B and D) ---
$ws.Range("B218").Value = 603042000000
$ws.Range("D218").Value = 131734713939.3145

$ws.Range("B219").Value = 603647900000
$ws.Range("D219").Value = 133394007027.1584

$ws.Range("B220").Value = 612183900000
$ws.Range("D220").Value = 131536473217.1634

# --- Append new row 224 (2023-06-01 data point) ---
# Copy formatting from the last existing data row (223) so the new date
# cell (A224) carries the same style (date number format / border / etc.)
$ws.Range("A223").Copy()
$ws.Range("A224").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("A224").Value = 45078
$ws.Range("B224").Value = 624519300000
$ws.Range("C224").Value = 0.2204342554833021
$ws.Range("D224").Value = 137665446930.453
